$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 17.38200411459175
$ws.Range("C2").Value = 10.24960728346503
$ws.Range("D2").Value = 14.96765855275356
$ws.Range("E2").Value = 16.38295555995843
$ws.Range("G2").Value = 3.664042355815601
$ws.Range("J2").Value = 9.351577344150392
$ws.Range("N2").Value = 17.63364990059747
$ws.Range("O2").Value = 27.44723639764635
$ws.Range("B3").Value = 16.88755965899315
$ws.Range("C3").Value = 9.794781807820081
$ws.Range("D3").Value = 14.90779738029639
$ws.Range("E3").Value = 16.32357327415852
$ws.Range("G3").Value = 3.66695021510899
$ws.Range("J3").Value = 9.361047719588052
$ws.Range("N3").Value = 17.69977067761061
$ws.Range("O3").Value = 27.46627507701007
$ws.Range("B4").Value = 16.57996326336046
$ws.Range("C4").Value = 9.505770003799208
$ws.Range("D4").Value = 14.87449218640738
$ws.Range("E4").Value = 16.29091949324159
$ws.Range("G4").Value = 3.668829198452866
$ws.Range("J4").Value = 9.368352338649997
$ws.Range("N4").Value = 17.74233256334899
$ws.Range("O4").Value = 27.48572465288974
$ws.Range("B5").Value = 16.4538065848335
$ws.Range("C5").Value = 9.385709511503421
$ws.Range("D5").Value = 14.86179682747303
$ws.Range("E5").Value = 16.27857884601919
$ws.Range("G5").Value = 3.669618505191222
$ws.Range("J5").Value = 9.371703537292724
$ws.Range("N5").Value = 17.76017213595967
$ws.Range("O5").Value = 27.49559575233886
$ws.Range("B6").Value = 16.43281555319487
$ws.Range("C6").Value = 9.365640772596279
$ws.Range("D6").Value = 14.85974197729597
$ws.Range("E6").Value = 16.27658826744258
$ws.Range("G6").Value = 3.669750997140372
$ws.Range("J6").Value = 9.372282615901719
$ws.Range("N6").Value = 17.76316434097231
$ws.Range("O6").Value = 27.49735213120061
$ws.Range("B7").Value = 16.57826487361607
$ws.Range("C7").Value = 9.504159846367227
$ws.Range("D7").Value = 14.87431741088485
$ws.Range("E7").Value = 16.29074914093289
$ws.Range("G7").Value = 3.668839747631286
$ws.Range("J7").Value = 9.368396017926328
$ws.Range("N7").Value = 17.74257114677303
$ws.Range("O7").Value = 27.48584991065317
$ws.Range("B8").Value = 17.21246647036027
$ws.Range("C8").Value = 10.09489869888206
$ws.Range("D8").Value = 14.94630832894101
$ws.Range("E8").Value = 16.36169593291265
$ws.Range("G8").Value = 3.665025620206651
$ws.Range("J8").Value = 9.35453348383874
$ws.Range("N8").Value = 17.65604174984166
$ws.Range("O8").Value = 27.45218710776402
$ws.Range("B9").Value = 18.41624751906374
$ws.Range("C9").Value = 11.16983042394258
$ws.Range("D9").Value = 15.1143956948594
$ws.Range("E9").Value = 16.53059607494911
$ws.Range("G9").Value = 3.658284593031968
$ws.Range("J9").Value = 9.339172528663056
$ws.Range("N9").Value = 17.5018669084232
$ws.Range("O9").Value = 27.44796567230115
$ws.Range("B10").Value = 19.26641713910166
$ws.Range("C10").Value = 11.90177865072323
$ws.Range("D10").Value = 15.25361085017808
$ws.Range("E10").Value = 16.67218172949586
$ws.Range("G10").Value = 3.653776858214577
$ws.Range("J10").Value = 9.335095454347838
$ws.Range("N10").Value = 17.39795076206115
$ws.Range("O10").Value = 27.48276515325113
$ws.Range("B11").Value = 19.64391413757073
$ws.Range("C11").Value = 12.22109160911201
$ws.Range("D11").Value = 15.32019004410622
$ws.Range("E11").Value = 16.74022922917665
$ws.Range("G11").Value = 3.651821650412472
$ws.Range("J11").Value = 9.334804473752351
$ws.Range("N11").Value = 17.35268737988725
$ws.Range("O11").Value = 27.50685033937085
$ws.Range("B12").Value = 19.78539802270125
$ws.Range("C12").Value = 12.33996941044319
$ws.Range("D12").Value = 15.34585359337868
$ws.Range("E12").Value = 16.76650512769332
$ws.Range("G12").Value = 3.651094893951785
$ws.Range("J12").Value = 9.334918834742643
$ws.Range("N12").Value = 17.3358346131796
$ws.Range("O12").Value = 27.51715755768818
$ws.Range("B13").Value = 19.75499416523218
$ws.Range("C13").Value = 12.31445871403573
$ws.Range("D13").Value = 15.34030667102722
$ws.Range("E13").Value = 16.76082381631134
$ws.Range("G13").Value = 3.651250808624515
$ws.Range("J13").Value = 9.334884224105622
$ws.Range("N13").Value = 17.33945139275013
$ws.Range("O13").Value = 27.51488495520113
$ws.Range("B14").Value = 19.65558421590453
$ws.Range("C14").Value = 12.23091301391664
$ws.Range("D14").Value = 15.32229244550961
$ws.Range("E14").Value = 16.74238087873083
$ws.Range("G14").Value = 3.651761586832476
$ws.Range("J14").Value = 9.334809384109011
$ws.Range("N14").Value = 17.35129513842055
$ws.Range("O14").Value = 27.50767454418614
$ws.Range("B15").Value = 19.59449801818765
$ws.Range("C15").Value = 12.17947127384505
$ws.Range("D15").Value = 15.31131652107026
$ws.Range("E15").Value = 16.73114969831069
$ws.Range("G15").Value = 3.652076227044057
$ws.Range("J15").Value = 9.334792774576533
$ws.Range("N15").Value = 17.35858717541564
$ws.Range("O15").Value = 27.50341246765193
$ws.Range("B16").Value = 19.24154877626026
$ws.Range("C16").Value = 11.88062929418561
$ws.Range("D16").Value = 15.24932379844403
$ws.Range("E16").Value = 16.66780659146491
$ws.Range("G16").Value = 3.653906548346401
$ws.Range("J16").Value = 9.335145911984466
$ws.Range("N16").Value = 17.40094913356361
$ws.Range("O16").Value = 27.48135723694165
$ws.Range("B17").Value = 19.02255590938442
$ws.Range("C17").Value = 11.69374448445209
$ws.Range("D17").Value = 15.21211456210982
$ws.Range("E17").Value = 16.62986902572484
$ws.Range("G17").Value = 3.655053765710897
$ws.Range("J17").Value = 9.335762861413519
$ws.Range("N17").Value = 17.42745031168444
$ws.Range("O17").Value = 27.46994147643008
$ws.Range("B18").Value = 18.89573368315222
$ws.Range("C18").Value = 11.58497178345837
$ws.Range("D18").Value = 15.19101965689671
$ws.Range("E18").Value = 16.6083918987987
$ws.Range("G18").Value = 3.655722597350956
$ws.Range("J18").Value = 9.336264921817177
$ws.Range("N18").Value = 17.44288222553967
$ws.Range("O18").Value = 27.46415251061293
$ws.Range("B19").Value = 18.85264993932381
$ws.Range("C19").Value = 11.54792562686892
$ws.Range("D19").Value = 15.18393045683039
$ws.Range("E19").Value = 16.60117957799621
$ws.Range("G19").Value = 3.655950597387954
$ws.Range("J19").Value = 9.336460201729187
$ws.Range("N19").Value = 17.44813973204639
$ws.Range("O19").Value = 27.46232591882142
$ws.Range("B20").Value = 19.04595845115325
$ws.Range("C20").Value = 11.71377190915281
$ws.Range("D20").Value = 15.21604390164189
$ws.Range("E20").Value = 16.6338720973156
$ws.Range("G20").Value = 3.654930713370595
$ws.Range("J20").Value = 9.335681952649029
$ws.Range("N20").Value = 17.42460965091037
$ws.Range("O20").Value = 27.47107627754616
$ws.Range("B21").Value = 19.68482413442192
$ws.Range("C21").Value = 12.25550829421391
$ws.Range("D21").Value = 15.32757153808364
$ws.Range("E21").Value = 16.74778436391118
$ws.Range("G21").Value = 3.651611189368559
$ws.Range("J21").Value = 9.334825275165525
$ws.Range("N21").Value = 17.34780855128684
$ws.Range("O21").Value = 27.50976021971193
$ws.Range("B22").Value = 20.09376064839457
$ws.Range("C22").Value = 12.59765213812335
$ws.Range("D22").Value = 15.40308480920108
$ws.Range("E22").Value = 16.82518454367096
$ws.Range("G22").Value = 3.649521147249028
$ws.Range("J22").Value = 9.335574020364575
$ws.Range("N22").Value = 17.29928972074494
$ws.Range("O22").Value = 27.54195821148472
$ws.Range("B23").Value = 19.87633192775171
$ws.Range("C23").Value = 12.41615542407803
$ws.Range("D23").Value = 15.36254733976624
$ws.Range("E23").Value = 16.78360988234324
$ws.Range("G23").Value = 3.650629397146214
$ws.Range("J23").Value = 9.335054787454716
$ws.Range("N23").Value = 17.32503229636335
$ws.Range("O23").Value = 27.52414119308541
$ws.Range("B24").Value = 19.03538101808908
$ws.Range("C24").Value = 11.70472164767679
$ws.Range("D24").Value = 15.21426652038794
$ws.Range("E24").Value = 16.6320612681174
$ws.Range("G24").Value = 3.654986316417726
$ws.Range("J24").Value = 9.335718072450824
$ws.Range("N24").Value = 17.42589330290458
$ws.Range("O24").Value = 27.47056082266676
$ws.Range("B25").Value = 18.09594073604566
$ws.Range("C25").Value = 10.88875384432977
$ws.Range("D25").Value = 15.06610900345792
$ws.Range("E25").Value = 16.4817799570047
$ws.Range("G25").Value = 3.66002971029148
$ws.Range("J25").Value = 9.34206187315079
$ws.Range("N25").Value = 17.5419252131518
$ws.Range("O25").Value = 27.49559575233886
